# Tracker updated And Python Task added
#
# - C15: "Control Flow tools, functions" -> add ",function args,"
# - C16: "tasks on lists,tuples,sets,dictionaries,methods etc." -> insert ",function args," before "etc."
# - New row 17: 2018-08-14, "pyhton" (python task entry), formatted like row 16's date/topic cells
# - Scroll the sheet down a bit (author had scrolled to show the new row)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new tracker row (row 17) with the same look as row 16 (date + topic styling).
$ws.Range("A16:B16").Copy()
$ws.Range("A17:B17").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("A17").Value = 43326
$ws.Range("B17").Value = "pyhton"

# Edit existing descriptions in column C for rows 15 and 16.
$ws.Range("C15").Value = "Control Flow tools, functions,function args,"
$ws.Range("C16").Value = "tasks on lists,tuples,sets,dictionaries,methods ,function args,etc."

# Restore the active cell/selection and nudge the visible scroll position down,
# matching the author having scrolled to the newly added row.
$ws.Range("C16").Select()
$win = $excel.ActiveWindow
$win.ScrollRow = 4
$win.ScrollColumn = 1
